# Re-shuffle of the "distr_elc_won" connector-table rows (cols S,T,AC,AD,AE,AG)
# on the "solar" sheet, rows 4-26. Each row keeps its left-hand block
# (columns B:O, the e_spv-CHE_00NN entries) fixed, but the right-hand
# connector block now pairs up with a different CHE_N grid cell than before,
# per the regenerated CHE_grids model data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$ws.Range("S4").Value = "distr_elc_won-CHE_0021"
$ws.Range("T4").Value = "connecting solar and wind to buses in grid cell CHE_21"
$ws.Range("AC4").Value = "distr_elc_won-CHE_0021"
$ws.Range("AD4").Value = "elc_won-CHE_0021,elc_spv-CHE_0021"
$ws.Range("AE4").Value = "e_CH17-380"
$ws.Range("AG4").Value = "CHE_21"
$ws.Range("S5").Value = "distr_elc_won-CHE_0009"
$ws.Range("T5").Value = "connecting solar and wind to buses in grid cell CHE_9"
$ws.Range("AC5").Value = "distr_elc_won-CHE_0009"
$ws.Range("AD5").Value = "elc_won-CHE_0009,elc_spv-CHE_0009"
$ws.Range("AE5").Value = "e_w1105061707-220,e_w165254212-220"
$ws.Range("AG5").Value = "CHE_9"
$ws.Range("S6").Value = "distr_elc_won-CHE_0004"
$ws.Range("T6").Value = "connecting solar and wind to buses in grid cell CHE_4"
$ws.Range("AC6").Value = "distr_elc_won-CHE_0004"
$ws.Range("AD6").Value = "elc_won-CHE_0004,elc_spv-CHE_0004"
$ws.Range("AE6").Value = "e_w194258388-220"
$ws.Range("AG6").Value = "CHE_4"
$ws.Range("S7").Value = "distr_elc_won-CHE_0017"
$ws.Range("T7").Value = "connecting solar and wind to buses in grid cell CHE_17"
$ws.Range("AC7").Value = "distr_elc_won-CHE_0017"
$ws.Range("AD7").Value = "elc_won-CHE_0017,elc_spv-CHE_0017"
$ws.Range("AE7").Value = "e_CH16-380,e_CH18-220"
$ws.Range("AG7").Value = "CHE_17"
$ws.Range("S8").Value = "distr_elc_won-CHE_0019"
$ws.Range("T8").Value = "connecting solar and wind to buses in grid cell CHE_19"
$ws.Range("AC8").Value = "distr_elc_won-CHE_0019"
$ws.Range("AD8").Value = "elc_won-CHE_0019,elc_spv-CHE_0019"
$ws.Range("AE8").Value = "e_CH45-220,e_w281804158-220"
$ws.Range("AG8").Value = "CHE_19"
$ws.Range("S9").Value = "distr_elc_won-CHE_0003"
$ws.Range("T9").Value = "connecting solar and wind to buses in grid cell CHE_3"
$ws.Range("AC9").Value = "distr_elc_won-CHE_0003"
$ws.Range("AD9").Value = "elc_won-CHE_0003,elc_spv-CHE_0003"
$ws.Range("AE9").Value = "e_CH48-225,e_CH49-225,e_CH60-225,e_w234983117-220,e_w260211728-225"
$ws.Range("AG9").Value = "CHE_3"
$ws.Range("S10").Value = "distr_elc_won-CHE_0024"
$ws.Range("T10").Value = "connecting solar and wind to buses in grid cell CHE_24"
$ws.Range("AC10").Value = "distr_elc_won-CHE_0024"
$ws.Range("AD10").Value = "elc_won-CHE_0024,elc_spv-CHE_0024"
$ws.Range("AE10").Value = "e_CH33-380,e_CH37-380,e_r7933294-380"
$ws.Range("AG10").Value = "CHE_24"
$ws.Range("S11").Value = "distr_elc_won-CHE_0005"
$ws.Range("T11").Value = "connecting solar and wind to buses in grid cell CHE_5"
$ws.Range("AC11").Value = "distr_elc_won-CHE_0005"
$ws.Range("AD11").Value = "elc_won-CHE_0005,elc_spv-CHE_0005"
$ws.Range("AE11").Value = "e_r5378910-220"
$ws.Range("AG11").Value = "CHE_5"
$ws.Range("S12").Value = "distr_elc_won-CHE_0008"
$ws.Range("T12").Value = "connecting solar and wind to buses in grid cell CHE_8"
$ws.Range("AC12").Value = "distr_elc_won-CHE_0008"
$ws.Range("AD12").Value = "elc_won-CHE_0008,elc_spv-CHE_0008"
$ws.Range("AE12").Value = "e_w209324991-220,e_w402055336-220,e_w758315582-220"
$ws.Range("AG12").Value = "CHE_8"
$ws.Range("S13").Value = "distr_elc_won-CHE_0014"
$ws.Range("T13").Value = "connecting solar and wind to buses in grid cell CHE_14"
$ws.Range("AC13").Value = "distr_elc_won-CHE_0014"
$ws.Range("AD13").Value = "elc_won-CHE_0014,elc_spv-CHE_0014"
$ws.Range("AE13").Value = "e_CH1-220,e_CH2-220,e_CH5-220,e_CH7-220,e_w108257952-220,e_w1284913429-220,e_w190819048-220,e_w33271433-220,e_w89405664-220,e_w89977424-220,e_w98648381-220"
$ws.Range("AG13").Value = "CHE_14"
$ws.Range("S14").Value = "distr_elc_won-CHE_0018"
$ws.Range("T14").Value = "connecting solar and wind to buses in grid cell CHE_18"
$ws.Range("AC14").Value = "distr_elc_won-CHE_0018"
$ws.Range("AD14").Value = "elc_won-CHE_0018,elc_spv-CHE_0018"
$ws.Range("AE14").Value = "e_CH6-220,e_w22899676-220,e_w240575085-220,e_w281809991-220,e_w50561341-220"
$ws.Range("AG14").Value = "CHE_18"
$ws.Range("S15").Value = "distr_elc_won-CHE_0012"
$ws.Range("T15").Value = "connecting solar and wind to buses in grid cell CHE_12"
$ws.Range("AC15").Value = "distr_elc_won-CHE_0012"
$ws.Range("AD15").Value = "elc_won-CHE_0012,elc_spv-CHE_0012"
$ws.Range("AE15").Value = "e_CH19-220,e_CH20-220,e_CH25-220,e_CH27-220,e_CH35-220,e_CH36-220,e_CH38-220,e_CH39-220,e_CH40-220,e_w1086214433-220,e_w240959264-220"
$ws.Range("AG15").Value = "CHE_12"
$ws.Range("S16").Value = "distr_elc_won-CHE_0000"
$ws.Range("T16").Value = "connecting solar and wind to buses in grid cell CHE_0"
$ws.Range("AC16").Value = "distr_elc_won-CHE_0000"
$ws.Range("AD16").Value = "elc_won-CHE_0000,elc_spv-CHE_0000"
$ws.Range("AE16").Value = "e_CH57-220"
$ws.Range("AG16").Value = "CHE_0"
$ws.Range("S17").Value = "distr_elc_won-CHE_0011"
$ws.Range("T17").Value = "connecting solar and wind to buses in grid cell CHE_11"
$ws.Range("AC17").Value = "distr_elc_won-CHE_0011"
$ws.Range("AD17").Value = "elc_won-CHE_0011,elc_spv-CHE_0011"
$ws.Range("AE17").Value = "e_CH46-220,e_w391576135-220,e_w969819301-220"
$ws.Range("AG17").Value = "CHE_11"
$ws.Range("S18").Value = "distr_elc_won-CHE_0015"
$ws.Range("T18").Value = "connecting solar and wind to buses in grid cell CHE_15"
$ws.Range("AC18").Value = "distr_elc_won-CHE_0015"
$ws.Range("AD18").Value = "elc_won-CHE_0015,elc_spv-CHE_0015"
$ws.Range("AE18").Value = "e_CH50-220,e_CH59-220,e_w1327084723-220,e_w281803398-220"
$ws.Range("AG18").Value = "CHE_15"
$ws.Range("S19").Value = "distr_elc_won-CHE_0025"
$ws.Range("T19").Value = "connecting solar and wind to buses in grid cell CHE_25"
$ws.Range("AC19").Value = "distr_elc_won-CHE_0025"
$ws.Range("AD19").Value = "elc_won-CHE_0025,elc_spv-CHE_0025"
$ws.Range("AE19").Value = "e_r7933294-380"
$ws.Range("AG19").Value = "CHE_25"
$ws.Range("S20").Value = "distr_elc_won-CHE_0013"
$ws.Range("T20").Value = "connecting solar and wind to buses in grid cell CHE_13"
$ws.Range("AC20").Value = "distr_elc_won-CHE_0013"
$ws.Range("AD20").Value = "elc_won-CHE_0013,elc_spv-CHE_0013"
$ws.Range("AE20").Value = "e_w146225999-220,e_w159527493-220,e_w281799252-220,e_w97941869-220"
$ws.Range("AG20").Value = "CHE_13"
$ws.Range("S21").Value = "distr_elc_won-CHE_0007"
$ws.Range("T21").Value = "connecting solar and wind to buses in grid cell CHE_7"
$ws.Range("AC21").Value = "distr_elc_won-CHE_0007"
$ws.Range("AD21").Value = "elc_won-CHE_0007,elc_spv-CHE_0007"
$ws.Range("AE21").Value = "e_CH44-220,e_w758943072-220"
$ws.Range("AG21").Value = "CHE_7"
$ws.Range("S22").Value = "distr_elc_won-CHE_0010"
$ws.Range("T22").Value = "connecting solar and wind to buses in grid cell CHE_10"
$ws.Range("AC22").Value = "distr_elc_won-CHE_0010"
$ws.Range("AD22").Value = "elc_won-CHE_0010,elc_spv-CHE_0010"
$ws.Range("AE22").Value = "e_CH11-220,e_w109037817-220,e_w127004407-380,e_w27435934-220,e_w30350721-220"
$ws.Range("AG22").Value = "CHE_10"
$ws.Range("S23").Value = "distr_elc_won-CHE_0022"
$ws.Range("T23").Value = "connecting solar and wind to buses in grid cell CHE_22"
$ws.Range("AC23").Value = "distr_elc_won-CHE_0022"
$ws.Range("AD23").Value = "elc_won-CHE_0022,elc_spv-CHE_0022"
$ws.Range("AE23").Value = "e_CH4-220,e_w455120191-220,e_w83861269-220"
$ws.Range("AG23").Value = "CHE_22"
$ws.Range("S24").Value = "distr_elc_won-CHE_0020"
$ws.Range("T24").Value = "connecting solar and wind to buses in grid cell CHE_20"
$ws.Range("AC24").Value = "distr_elc_won-CHE_0020"
$ws.Range("AD24").Value = "elc_won-CHE_0020,elc_spv-CHE_0020"
$ws.Range("AE24").Value = "e_CH21-220,e_CH22-220,e_CH29-220,e_CH41-380,e_w1208713169-220,e_w207993342-220,e_w208780268-380,e_w212498548-220,e_w36348118-220,e_w71500123-220"
$ws.Range("AG24").Value = "CHE_20"
$ws.Range("S25").Value = "distr_elc_won-CHE_0001"
$ws.Range("T25").Value = "connecting solar and wind to buses in grid cell CHE_1"
$ws.Range("AC25").Value = "distr_elc_won-CHE_0001"
$ws.Range("AD25").Value = "elc_won-CHE_0001,elc_spv-CHE_0001"
$ws.Range("AE25").Value = "e_CH31-220,e_w132373704-220,e_w55695765-220"
$ws.Range("AG25").Value = "CHE_1"
$ws.Range("S26").Value = "distr_elc_won-CHE_0006"
$ws.Range("T26").Value = "connecting solar and wind to buses in grid cell CHE_6"
$ws.Range("AC26").Value = "distr_elc_won-CHE_0006"
$ws.Range("AD26").Value = "elc_won-CHE_0006,elc_spv-CHE_0006"
$ws.Range("AE26").Value = "e_w127004407-380"
$ws.Range("AG26").Value = "CHE_6"
